$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card6")

# O1 header: remove trailing space "Serviced by " -> "Serviced by"
$ws.Range("O1").Value = "Serviced by"

# O2: was an empty inline string cell, becomes "nan"
$ws.Range("O2").Value = "nan"

# Row 3 edits: M3, N3 get real values, O3 gets a value (was empty before)
$ws.Range("M3").Value = "سيرفيس"
$ws.Range("N3").Value = "تم تغير اول جريده خلفي20وسن فلاتس المتحرك"
$ws.Range("O3").Value = "م.محمد عبدالله ،محمود ايهاب"

# O4..O12: were empty inline string cells, become "nan"
for ($r = 4; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = "nan"
}
